# Rename the worksheet "Uncut_Sheet" to "Uncut_Sheet_1".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Uncut_Sheet")
$ws.Name = "Uncut_Sheet_1"

# The workbook-scoped Print_Area defined name still points at the old
# sheet name after the rename, so repoint it explicitly.
for ($i = 1; $i -le $wb.Names.Count; $i++) {
    $nm = $wb.Names.Item($i)
    if ($nm.Name -like "*Print_Area*") {
        $nm.RefersTo = "=" + $ws.Name + "!`$A`$1:`$G`$42"
    }
}

# Move the active selection to B21:C21.
$ws.Range("B21:C21").Select()
